$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the higher-numbered row first so the other row's index stays valid.
# Row 28 = "SC 92"
$ws.Rows.Item(28).Delete()
# Row 26 = "RM 232"
$ws.Rows.Item(26).Delete()
